$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Append a new data row (row 25) to the bottom of the table.
$ws.Range("A25").Value = "21/06/2024 05:44:39"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "TATATECH"
$ws.Range("D25").Value = "Tata Technologies Ltd"

# bsecode is stored as text (not a number) for this row in the source data,
# so force the cell to text before writing the digit string, then restore
# the default (Normal) style so no stray formatting is left behind.
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "544028"
$ws.Range("E25").Style = "Normal"

$ws.Range("F25").Value = -0.8100000000000001
$ws.Range("G25").Value = 1002.15
$ws.Range("H25").Value = 1447221
